$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).NumberFormat = '@'
$ws.Cells.Item(2,4).Value = '63.064.96'
$ws.Cells.Item(2,4).Style = 'Normal'
$ws.Cells.Item(2,5).Value = '  -4.43%  '
$ws.Cells.Item(3,4).NumberFormat = '@'
$ws.Cells.Item(3,4).Value = '3.089.37'
$ws.Cells.Item(3,4).Style = 'Normal'
$ws.Cells.Item(3,5).Value = '  -4.56%  '
$ws.Cells.Item(4,4).NumberFormat = '@'
$ws.Cells.Item(4,4).Value = '0.998'
$ws.Cells.Item(4,4).Style = 'Normal'
$ws.Cells.Item(4,5).Value = '  -0.28%  '
$ws.Cells.Item(5,4).NumberFormat = '@'
$ws.Cells.Item(5,4).Value = '542.40'
$ws.Cells.Item(5,4).Style = 'Normal'
$ws.Cells.Item(5,5).Value = '  -6.47%  '
$ws.Cells.Item(6,4).NumberFormat = '@'
$ws.Cells.Item(6,4).Value = '135.61'
$ws.Cells.Item(6,4).Style = 'Normal'
$ws.Cells.Item(6,5).Value = '  -10.34%  '
$ws.Cells.Item(7,5).Value = '  -0.02%  '
$ws.Cells.Item(8,4).NumberFormat = '@'
$ws.Cells.Item(8,4).Value = '3.086.98'
$ws.Cells.Item(8,4).Style = 'Normal'
$ws.Cells.Item(8,5).Value = '  -4.34%  '
$ws.Cells.Item(9,4).NumberFormat = '@'
$ws.Cells.Item(9,4).Value = '0.492'
$ws.Cells.Item(9,4).Style = 'Normal'
$ws.Cells.Item(9,5).Value = '  -4.01%  '
$ws.Cells.Item(10,4).NumberFormat = '@'
$ws.Cells.Item(10,4).Value = '0.156'
$ws.Cells.Item(10,4).Style = 'Normal'
$ws.Cells.Item(10,5).Value = '  -4.14%  '
$ws.Cells.Item(11,4).NumberFormat = '@'
$ws.Cells.Item(11,4).Value = '6.26'
$ws.Cells.Item(11,4).Style = 'Normal'
$ws.Cells.Item(11,5).Value = '  -11.37%  '
$ws.Cells.Item(12,4).NumberFormat = '@'
$ws.Cells.Item(12,4).Value = '0.465'
$ws.Cells.Item(12,4).Style = 'Normal'
$ws.Cells.Item(12,5).Value = '  -4.38%  '
$ws.Cells.Item(13,4).NumberFormat = '@'
$ws.Cells.Item(13,4).Value = '35.11'
$ws.Cells.Item(13,4).Style = 'Normal'
$ws.Cells.Item(13,5).Value = '  -6.54%  '
$ws.Cells.Item(14,4).NumberFormat = '@'
$ws.Cells.Item(14,4).Value = '0.0000223'
$ws.Cells.Item(14,4).Style = 'Normal'
$ws.Cells.Item(14,5).Value = '  -4.12%  '
$ws.Cells.Item(15,4).NumberFormat = '@'
$ws.Cells.Item(15,4).Value = '3.553.32'
$ws.Cells.Item(15,4).Style = 'Normal'
$ws.Cells.Item(15,5).Value = '  -5.29%  '
$ws.Cells.Item(16,4).NumberFormat = '@'
$ws.Cells.Item(16,4).Value = '62.884.38'
$ws.Cells.Item(16,4).Style = 'Normal'
$ws.Cells.Item(16,5).Value = '  -4.83%  '
$ws.Cells.Item(17,5).Value = '  -2.81%  '
$ws.Cells.Item(18,4).NumberFormat = '@'
$ws.Cells.Item(18,4).Value = '3.078.23'
$ws.Cells.Item(18,4).Style = 'Normal'
$ws.Cells.Item(18,5).Value = '  -5.03%  '
$ws.Cells.Item(19,4).NumberFormat = '@'
$ws.Cells.Item(19,4).Value = '6.69'
$ws.Cells.Item(19,4).Style = 'Normal'
$ws.Cells.Item(19,5).Value = '  -5.56%  '
$ws.Cells.Item(20,4).NumberFormat = '@'
$ws.Cells.Item(20,4).Value = '485.79'
$ws.Cells.Item(20,4).Style = 'Normal'
$ws.Cells.Item(20,5).Value = '  -10.66%  '
$ws.Cells.Item(21,4).NumberFormat = '@'
$ws.Cells.Item(21,4).Value = '13.49'
$ws.Cells.Item(21,4).Style = 'Normal'
$ws.Cells.Item(21,5).Value = '  -6.76%  '
$ws.Cells.Item(22,4).NumberFormat = '@'
$ws.Cells.Item(22,4).Value = '0.711'
$ws.Cells.Item(22,4).Style = 'Normal'
$ws.Cells.Item(22,5).Value = '  -4.28%  '
$ws.Cells.Item(23,4).NumberFormat = '@'
$ws.Cells.Item(23,4).Value = '7.29'
$ws.Cells.Item(23,4).Style = 'Normal'
$ws.Cells.Item(23,5).Value = '  -7.18%  '
$ws.Cells.Item(24,4).NumberFormat = '@'
$ws.Cells.Item(24,4).Value = '79.00'
$ws.Cells.Item(24,4).Style = 'Normal'
$ws.Cells.Item(24,5).Value = '  -2.23%  '
$ws.Cells.Item(25,4).NumberFormat = '@'
$ws.Cells.Item(25,4).Value = '12.17'
$ws.Cells.Item(25,4).Style = 'Normal'
$ws.Cells.Item(25,5).Value = '  -9.36%  '
$ws.Cells.Item(26,5).Value = '  +0.20%  '
$ws.Cells.Item(27,4).NumberFormat = '@'
$ws.Cells.Item(27,4).Value = '2.73'
$ws.Cells.Item(27,4).Style = 'Normal'
$ws.Cells.Item(27,5).Value = '  -8.12%  '
$ws.Cells.Item(28,4).NumberFormat = '@'
$ws.Cells.Item(28,4).Value = '8.29'
$ws.Cells.Item(28,4).Style = 'Normal'
$ws.Cells.Item(28,5).Value = '  -10.70%  '
$ws.Cells.Item(29,4).NumberFormat = '@'
$ws.Cells.Item(29,4).Value = '0.996'
$ws.Cells.Item(29,4).Style = 'Normal'
$ws.Cells.Item(29,5).Value = '  -0.42%  '
$ws.Cells.Item(30,4).NumberFormat = '@'
$ws.Cells.Item(30,4).Value = '1.94'
$ws.Cells.Item(30,4).Style = 'Normal'
$ws.Cells.Item(30,5).Value = '  -13.44%  '
$ws.Cells.Item(31,4).NumberFormat = '@'
$ws.Cells.Item(31,4).Value = '26.22'
$ws.Cells.Item(31,4).Style = 'Normal'
$ws.Cells.Item(31,5).Value = '  -4.95%  '
$ws.Cells.Item(32,4).NumberFormat = '@'
$ws.Cells.Item(32,4).Value = '1.11'
$ws.Cells.Item(32,4).Style = 'Normal'
$ws.Cells.Item(32,5).Value = '  -5.90%  '
$ws.Cells.Item(33,2).Value = 'Stacks'
$ws.Cells.Item(33,3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(33,4).NumberFormat = '@'
$ws.Cells.Item(33,4).Value = '2.44'
$ws.Cells.Item(33,4).Style = 'Normal'
$ws.Cells.Item(33,5).Value = '  -10.83%  '
$ws.Cells.Item(34,2).Value = 'OKB'
$ws.Cells.Item(34,3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(34,4).NumberFormat = '@'
$ws.Cells.Item(34,4).Value = '59.07'
$ws.Cells.Item(34,4).Style = 'Normal'
$ws.Cells.Item(34,5).Value = '  +7.12%  '
$ws.Cells.Item(35,4).NumberFormat = '@'
$ws.Cells.Item(35,4).Value = '6.06'
$ws.Cells.Item(35,4).Style = 'Normal'
$ws.Cells.Item(35,5).Value = '  -4.18%  '
$ws.Cells.Item(36,2).Value = 'NEARProtocol'
$ws.Cells.Item(36,3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(36,4).NumberFormat = '@'
$ws.Cells.Item(36,4).Value = '5.23'
$ws.Cells.Item(36,4).Style = 'Normal'
$ws.Cells.Item(36,5).Value = '  -7.11%  '
$ws.Cells.Item(37,2).Value = 'Bittensor'
$ws.Cells.Item(37,3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(37,4).NumberFormat = '@'
$ws.Cells.Item(37,4).Value = '480.49'
$ws.Cells.Item(37,4).Style = 'Normal'
$ws.Cells.Item(37,5).Value = '  -15.56%  '
$ws.Cells.Item(38,4).NumberFormat = '@'
$ws.Cells.Item(38,4).Value = '3.155.99'
$ws.Cells.Item(38,4).Style = 'Normal'
$ws.Cells.Item(38,5).Value = '  -1.17%  '
$ws.Cells.Item(39,2).Value = 'VeChain'
$ws.Cells.Item(39,3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(39,4).NumberFormat = '@'
$ws.Cells.Item(39,4).Value = '0.0394'
$ws.Cells.Item(39,4).Style = 'Normal'
$ws.Cells.Item(39,5).Value = '  -12.95%  '
$ws.Cells.Item(40,2).Value = 'Hedera'
$ws.Cells.Item(40,3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(40,4).NumberFormat = '@'
$ws.Cells.Item(40,4).Value = '0.0802'
$ws.Cells.Item(40,4).Style = 'Normal'
$ws.Cells.Item(40,5).Value = '  -6.41%  '
$ws.Cells.Item(41,4).NumberFormat = '@'
$ws.Cells.Item(41,4).Value = '0.117'
$ws.Cells.Item(41,4).Style = 'Normal'
$ws.Cells.Item(41,5).Value = '  -9.78%  '
$ws.Cells.Item(42,4).NumberFormat = '@'
$ws.Cells.Item(42,4).Value = '8.12'
$ws.Cells.Item(42,4).Style = 'Normal'
$ws.Cells.Item(42,5).Value = '  -5.08%  '
$ws.Cells.Item(43,4).NumberFormat = '@'
$ws.Cells.Item(43,4).Value = '2.59'
$ws.Cells.Item(43,4).Style = 'Normal'
$ws.Cells.Item(43,5).Value = '  -11.87%  '
$ws.Cells.Item(44,4).NumberFormat = '@'
$ws.Cells.Item(44,4).Value = '0.255'
$ws.Cells.Item(44,4).Style = 'Normal'
$ws.Cells.Item(44,5).Value = '  -9.13%  '
$ws.Cells.Item(45,5).Value = '  +0.05%  '
$ws.Cells.Item(46,4).NumberFormat = '@'
$ws.Cells.Item(46,4).Value = '2.07'
$ws.Cells.Item(46,4).Style = 'Normal'
$ws.Cells.Item(46,5).Value = '  -9.64%  '
$ws.Cells.Item(47,4).NumberFormat = '@'
$ws.Cells.Item(47,4).Value = '24.95'
$ws.Cells.Item(47,4).Style = 'Normal'
$ws.Cells.Item(47,5).Value = '  -5.31%  '
$ws.Cells.Item(48,4).NumberFormat = '@'
$ws.Cells.Item(48,4).Value = '118.99'
$ws.Cells.Item(48,4).Style = 'Normal'
$ws.Cells.Item(48,5).Value = '  -5.55%  '
$ws.Cells.Item(49,4).NumberFormat = '@'
$ws.Cells.Item(49,4).Value = '0.108'
$ws.Cells.Item(49,4).Style = 'Normal'
$ws.Cells.Item(49,5).Value = '  -3.47%  '
$ws.Cells.Item(50,4).NumberFormat = '@'
$ws.Cells.Item(50,4).Value = '0.0₃0515'
$ws.Cells.Item(50,4).Style = 'Normal'
$ws.Cells.Item(50,5).Value = '  -7.43%  '
$ws.Cells.Item(51,4).NumberFormat = '@'
$ws.Cells.Item(51,4).Value = '2.03'
$ws.Cells.Item(51,4).Style = 'Normal'
$ws.Cells.Item(51,5).Value = '  -7.72%  '
